# Apply updated crypto market data (prices + volume deltas) per upstream refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.687.30"
$ws.Range("E2").Value = "  +2.32%  "
$ws.Range("D3").Value = "2.526.84"
$ws.Range("E3").Value = "  +2.36%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Formula = "'594.10"
$ws.Range("E5").Value = "  +2.01%  "
$ws.Range("D6").Formula = "'176.99"
$ws.Range("E6").Value = "  +1.66%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +1.55%  "
$ws.Range("D9").Value = "2.526.58"
$ws.Range("E9").Value = "  +2.36%  "
$ws.Range("E10").Value = "  +5.94%  "
$ws.Range("E11").Value = "  -1.13%  "
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("E13").Value = "  +1.59%  "
$ws.Range("D14").Value = "2.989.07"
$ws.Range("E14").Value = "  +2.72%  "
$ws.Range("D15").Formula = "'26.23"
$ws.Range("E15").Value = "  +3.37%  "
$ws.Range("D16").Value = "68.573.36"
$ws.Range("E16").Value = "  +2.23%  "
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").Value = "2.533.74"
$ws.Range("E18").Value = "  +2.64%  "
$ws.Range("E19").Value = "  +1.72%  "
$ws.Range("D20").Formula = "'7.50"
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("D21").Formula = "'352.59"
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("D22").Formula = "'4.20"
$ws.Range("E22").Value = "  +4.66%  "
$ws.Range("D23").Formula = "'1.00"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Formula = "'70.85"
$ws.Range("E24").Value = "  +2.01%  "
$ws.Range("D25").Formula = "'4.24"
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("E26").Value = "  -5.29%  "
$ws.Range("E27").Value = "  -2.29%  "
$ws.Range("D28").Value = "2.690.84"
$ws.Range("D29").Formula = "'0.995"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").Value = "0.0₃0894"
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("D31").Formula = "'508.33"
$ws.Range("E31").Value = "  +1.87%  "
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("E33").Value = "  +1.82%  "
$ws.Range("E34").Value = "  +1.24%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").Formula = "'162.96"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("D39").Formula = "'18.41"
$ws.Range("E39").Value = "  +1.36%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Formula = "'1.77"
$ws.Range("E40").Value = "  +5.49%  "
$ws.Range("B41").Value = "ImmutableX"
$ws.Range("C41").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D41").Formula = "'1.32"
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").Formula = "'4.86"
$ws.Range("E43").Value = "  +1.13%  "
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("E45").Value = "  +1.46%  "
$ws.Range("D46").Formula = "'152.92"
$ws.Range("E46").Value = "  +7.18%  "
$ws.Range("E47").Value = "  +2.59%  "
$ws.Range("E48").Value = "  +2.52%  "
$ws.Range("E49").Value = "  +1.38%  "
$ws.Range("E50").Value = "  +2.41%  "
$ws.Range("E51").Value = "  -0.16%  "
